$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column stays as text so values like "1.00" are not
# coerced into numbers by Excel's automatic type detection.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '62.640.56'
$ws.Range("E2").Value = '  -1.52%  '
$ws.Range("D3").Value = '3.023.39'
$ws.Range("E3").Value = '  -1.69%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").Value = '585.87'
$ws.Range("E5").Value = '  -1.00%  '
$ws.Range("D6").Value = '148.14'
$ws.Range("E6").Value = '  -4.16%  '
$ws.Range("E7").Value = '  +0.23%  '
$ws.Range("E8").Value = '  -2.15%  '
$ws.Range("D9").Value = '3.018.48'
$ws.Range("E9").Value = '  -2.03%  '
$ws.Range("D10").Value = '0.150'
$ws.Range("E10").Value = '  -4.22%  '
$ws.Range("D11").Value = '5.86'
$ws.Range("E11").Value = '  -0.89%  '
$ws.Range("D12").Value = '0.457'
$ws.Range("E12").Value = '  +0.92%  '
$ws.Range("D13").Value = '0.0000230'
$ws.Range("E13").Value = '  -3.42%  '
$ws.Range("D14").Value = '34.81'
$ws.Range("E14").Value = '  -5.61%  '
$ws.Range("E15").Value = '  +1.85%  '
$ws.Range("D16").Value = '3.522.17'
$ws.Range("E16").Value = '  -1.61%  '
$ws.Range("E17").Value = '  -0.84%  '
$ws.Range("D18").Value = '62.619.32'
$ws.Range("E18").Value = '  -1.37%  '
$ws.Range("D19").Value = '3.019.97'
$ws.Range("E19").Value = '  -1.52%  '
$ws.Range("D20").Value = '465.07'
$ws.Range("E20").Value = '  -4.14%  '
$ws.Range("E21").Value = '  -3.84%  '
$ws.Range("D22").Value = '0.689'
$ws.Range("E22").Value = '  -3.02%  '
$ws.Range("D23").Value = '7.49'
$ws.Range("E23").Value = '  -1.49%  '
$ws.Range("D24").Value = '81.77'
$ws.Range("E24").Value = '  -0.35%  '
$ws.Range("D25").Value = '2.27'
$ws.Range("E25").Value = '  -6.31%  '
$ws.Range("D26").Value = '12.42'
$ws.Range("E26").Value = '  -4.28%  '
$ws.Range("D27").Value = '10.37'
$ws.Range("E27").Value = '  -1.73%  '
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.17%  '
$ws.Range("E30").Value = '  -1.95%  '
$ws.Range("D31").Value = '7.14'
$ws.Range("E31").Value = '  -5.24%  '
$ws.Range("E32").Value = '  -5.53%  '
$ws.Range("D33").Value = '28.82'
$ws.Range("E33").Value = '  +5.23%  '
$ws.Range("E34").Value = '  -2.92%  '
$ws.Range("D35").Value = '0.0₃0809'
$ws.Range("E35").Value = '  -1.95%  '
$ws.Range("E36").Value = '  -4.16%  '
$ws.Range("D37").Value = '5.80'
$ws.Range("E37").Value = '  -4.47%  '
$ws.Range("D38").Value = '2.14'
$ws.Range("E38").Value = '  -4.12%  '
$ws.Range("D39").Value = '50.49'
$ws.Range("E39").Value = '  -0.40%  '
$ws.Range("E40").Value = '  -1.89%  '
$ws.Range("D41").Value = '2.96'
$ws.Range("E41").Value = '  -8.90%  '
$ws.Range("D42").Value = '0.115'
$ws.Range("E42").Value = '  +2.17%  '
$ws.Range("D43").Value = '399.65'
$ws.Range("E43").Value = '  -9.35%  '
$ws.Range("D44").Value = '0.278'
$ws.Range("E44").Value = '  -4.19%  '
$ws.Range("D45").Value = '0.0360'
$ws.Range("E45").Value = '  -1.27%  '
$ws.Range("D46").Value = '2.756.11'
$ws.Range("E46").Value = '  -2.73%  '
$ws.Range("D47").Value = '37.40'
$ws.Range("E47").Value = '  -5.89%  '
$ws.Range("D48").Value = '129.01'
$ws.Range("E48").Value = '  -2.95%  '
$ws.Range("E49").Value = '  +0.08%  '
$ws.Range("E50").Value = '  -0.44%  '
$ws.Range("E51").Value = '  -1.62%  '
